# Update "Pais" worksheet with refreshed COVID-19 figures and new timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp in A1 (shared string used by the header row)
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 08:52"

# Updated per-country counters (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 37 - Chequia
$ws.Range("B37").Value = 6657
$ws.Range("C37").Value = 51
$ws.Range("D37").Value = 1235
$ws.Range("E37").Value = 5241
$ws.Range("F37").Value = 84

# Row 44 - Malasia
$ws.Range("B44").Value = 5449
$ws.Range("C44").Value = 343
$ws.Range("D44").Value = 347
$ws.Range("E44").Value = 4961
$ws.Range("F44").Value = 45
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 141

# Row 45 - Ucrania
$ws.Range("B45").Value = 5305
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 3102
$ws.Range("E45").Value = 2115
$ws.Range("F45").Value = 49
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 88

# Row 62 - Hungria
$ws.Range("B62").Value = 1916
$ws.Range("C62").Value = 82
$ws.Range("D62").Value = 250
$ws.Range("E62").Value = 1494
$ws.Range("F62").Value = 61

# Row 74 - Bosnia y Herzegovina
$ws.Range("B74").Value = 1298
$ws.Range("C74").Value = 59
$ws.Range("D74").Value = 242
$ws.Range("E74").Value = 1023
$ws.Range("F74").Value = 14
$ws.Range("H74").Value = 33

# Row 75 - Armenia
$ws.Range("B75").Value = 1268
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 338
$ws.Range("E75").Value = 883
$ws.Range("F75").Value = 4
$ws.Range("H75").Value = 47

# Row 76 - Lituania
$ws.Range("B76").Value = 1266
$ws.Range("C76").Value = 86
$ws.Range("D76").Value = 233
$ws.Range("E76").Value = 1027
$ws.Range("F76").Value = 3
$ws.Range("H76").Value = 6

# Row 77 - Oman
$ws.Range("B77").Value = 1248
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 523
$ws.Range("E77").Value = 705
$ws.Range("F77").Value = 30
$ws.Range("H77").Value = 20

# Row 106 - Estado de Palestina
$ws.Range("B106").Value = 420
$ws.Range("C106").Value = 22
$ws.Range("D106").Value = 189
$ws.Range("E106").Value = 225
$ws.Range("H106").Value = 6

# Row 107 - Jordania
$ws.Range("B107").Value = 418
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 69
$ws.Range("E107").Value = 347
$ws.Range("F107").Value = 0
$ws.Range("H107").Value = 2

# Row 108 - Reunion
$ws.Range("B108").Value = 413
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 269
$ws.Range("E108").Value = 137
$ws.Range("F108").Value = 5
$ws.Range("H108").Value = 7

# Row 109 - Taiwan
$ws.Range("B109").Value = 407
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 237
$ws.Range("E109").Value = 170
$ws.Range("F109").Value = 4
$ws.Range("H109").Value = 0

# Row 113 - Consejo Danes para los Refugiados
$ws.Range("B113").Value = 308
$ws.Range("C113").Value = 1
$ws.Range("D113").Value = 55
$ws.Range("E113").Value = 248
$ws.Range("F113").Value = 7
$ws.Range("H113").Value = 5

# Row 114 - Montenegro
$ws.Range("D114").Value = 26
$ws.Range("E114").Value = 256
$ws.Range("F114").Value = 0
$ws.Range("H114").Value = 25
